$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-39
# from serial date 45203 (2023-10-04) to 45205 (2023-10-06)
for ($row = 2; $row -le 39; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
